# Apply the "introduced guest account and handling" edit.
#
# users sheet:
#   - rick's count (C7) goes from 28 to 49
#   - a new guest/gast row is appended at row 8 with count 13
#   - selection moves from B8 to A8
#
# actions sheet:
#   - 34 new action-log rows are appended (rows 16-49)

$wb = $excel.ActiveWorkbook

# --- users sheet -----------------------------------------------------
$wsUsers = $wb.Worksheets.Item("users")

$wsUsers.Cells.Item(7, 3).Value = 49

$wsUsers.Cells.Item(8, 1).Value = "guest"
$wsUsers.Cells.Item(8, 2).Value = "gast"
$wsUsers.Cells.Item(8, 3).Value = 13

$wsUsers.Activate()
$wsUsers.Range("A8").Select()

# --- actions sheet -----------------------------------------------------
$wsActions = $wb.Worksheets.Item("actions")

$newActions = @(
  @("4ot355g53h", "rick", "2025-09-11T16:37:34.569542"),
  @("4ot355g53h", "rick", "2025-09-11T16:45:03.859558"),
  @("4ot355g53h", "rick", "2025-09-11T16:45:04.634407"),
  @("guest", "gast", "2025-09-11T16:46:44.827159"),
  @("guest", "gast", "2025-09-11T16:46:48.256704"),
  @("guest", "gast", "2025-09-11T16:46:49.158282"),
  @("guest", "gast", "2025-09-11T16:46:50.935001"),
  @("4ot355g53h", "rick", "2025-09-11T16:47:04.566519"),
  @("4ot355g53h", "rick", "2025-09-11T16:47:06.544106"),
  @("4ot355g53h", "rick", "2025-09-11T16:48:57.691330"),
  @("4ot355g53h", "rick", "2025-09-11T17:01:16.929735"),
  @("4ot355g53h", "rick", "2025-09-11T17:01:19.997282"),
  @("4ot355g53h", "rick", "2025-09-11T17:01:21.817021"),
  @("4ot355g53h", "rick", "2025-09-11T17:01:24.021084"),
  @("4ot355g53h", "rick", "2025-09-11T17:01:24.335717"),
  @("4ot355g53h", "rick", "2025-09-11T17:01:24.530092"),
  @("4ot355g53h", "rick", "2025-09-11T17:01:26.602666"),
  @("4ot355g53h", "rick", "2025-09-11T17:01:26.774413"),
  @("4ot355g53h", "rick", "2025-09-11T17:01:26.946217"),
  @("4ot355g53h", "rick", "2025-09-11T17:01:28.218388"),
  @("4ot355g53h", "rick", "2025-09-11T17:10:08.959011"),
  @("4ot355g53h", "rick", "2025-09-11T17:10:09.711290"),
  @("guest", "gast", "2025-09-11T17:10:11.195980"),
  @("guest", "gast", "2025-09-11T17:10:12.205959"),
  @("guest", "gast", "2025-09-11T17:10:12.782371"),
  @("guest", "gast", "2025-09-11T17:10:14.044459"),
  @("4ot355g53h", "rick", "2025-09-11T17:10:15.672224"),
  @("4ot355g53h", "rick", "2025-09-11T17:10:15.877200"),
  @("guest", "gast", "2025-09-11T17:17:26.899329"),
  @("guest", "gast", "2025-09-11T17:23:14.373956"),
  @("guest", "gast", "2025-09-11T17:23:14.870656"),
  @("guest", "gast", "2025-09-11T17:23:15.082736"),
  @("guest", "gast", "2025-09-11T17:23:15.409804"),
  @("4ot355g53h", "rick", "2025-09-11T17:23:17.129687")
)

$r = 16
foreach ($row in $newActions) {
  $wsActions.Cells.Item($r, 1).Value = $row[0]
  $wsActions.Cells.Item($r, 2).Value = $row[1]
  $wsActions.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}
